$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for every existing data
# row (2..306) from 2023-10-03 (45202) to 2023-10-04 (45203).
for ($r = 2; $r -le 306; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# Row 306 gets an explicit row height (matches customHeight flag in target).
$ws.Rows.Item(306).RowHeight = 15

# Append the new record as row 307.
$ws.Range("A307").Value = "A 47288-2023"

$ws.Range("B307").Value = 45202
$ws.Range("B307").NumberFormat = "YYYY-MM-DD"

$ws.Range("C307").Value = 45203
$ws.Range("C307").NumberFormat = "YYYY-MM-DD"

$ws.Range("D307").Value = "GÄVLEBORGS LÄN"
$ws.Range("E307").Value = "SÖDERHAMN"

$ws.Range("G307").Value = 0.5
$ws.Range("H307").Value = 0
$ws.Range("I307").Value = 0
$ws.Range("J307").Value = 0
$ws.Range("K307").Value = 0
$ws.Range("L307").Value = 0
$ws.Range("M307").Value = 0
$ws.Range("N307").Value = 0
$ws.Range("O307").Value = 0
$ws.Range("P307").Value = 0
$ws.Range("Q307").Value = 0

$ws.Range("R307").Value = ""
$ws.Range("R307").WrapText = $true
